# Apply the changes described by the commit:
#  - Rename the " SVM Regression " label (row 7, col A) to "SVM Regression "
#    (drop the leading space). This causes the shared-string table to drop
#    the now-unused old string and append the new one at the end.
#  - Update the Random Forest Regression row's ILI_lagwk3 R2 train/test
#    scores (J6/K6).
#  - Touch cell J11 (copy formatting from the already-present "marker" cell
#    E14) so a new, otherwise-empty row 11 appears in the sheet, matching
#    the author's click-through edit.
#  - Leave the active selection on the SVM Regression row (row 7), as the
#    author had it selected when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the SVM Regression row label (remove the stray leading space).
$ws.Range("A7").Value = "SVM Regression "

# Update the ILI_lagwk3 R2 Train/Test scores for Random Forest Regression.
$ws.Range("J6").Value = 0.82699999999999996
$ws.Range("K6").Value = 0.47699999999999998

# Stamp J11 with the same formatting already used on E14 so a new row 11
# shows up in the sheet (empty value, inherited style).
$ws.Range("E14").Copy()
$ws.Range("J11").PasteSpecial(-4122)

# Leave the whole SVM Regression row selected, as in the saved file.
$ws.Rows("7:7").Select()
